$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.01780628407291
$ws.Range("C2").Value = 11.59622974356034
$ws.Range("D2").Value = 5.010288421096508
$ws.Range("F2").Value = 23.86193444520829
$ws.Range("G2").Value = 28.01021468317829
$ws.Range("H2").Value = 14.16238112059023
$ws.Range("K2").Value = 7.983569668974547
$ws.Range("L2").Value = 10.70363282952417
$ws.Range("N2").Value = 18.58677497725083
$ws.Range("O2").Value = 21.44664335933064
$ws.Range("B3").Value = 11.73860415932579
$ws.Range("C3").Value = 11.63632186442802
$ws.Range("D3").Value = 4.93762201859681
$ws.Range("F3").Value = 23.88924901470408
$ws.Range("G3").Value = 28.05597749221388
$ws.Range("H3").Value = 14.20302510177842
$ws.Range("K3").Value = 7.78295201998305
$ws.Range("L3").Value = 10.67585856233755
$ws.Range("N3").Value = 18.6420438613589
$ws.Range("O3").Value = 21.50898729329739
$ws.Range("B4").Value = 11.56586649067068
$ws.Range("C4").Value = 11.66219071995775
$ws.Range("D4").Value = 4.891768889208749
$ws.Range("F4").Value = 23.91185150398377
$ws.Range("G4").Value = 28.09228281088726
$ws.Range("H4").Value = 14.23000736311622
$ws.Range("K4").Value = 7.658094264360084
$ws.Range("L4").Value = 10.6609215325415
$ws.Range("N4").Value = 18.67760148224282
$ws.Range("O4").Value = 21.55142090652961
$ws.Range("B5").Value = 11.49525063322107
$ws.Range("C5").Value = 11.67304826482828
$ws.Range("D5").Value = 4.8727832258032
$ws.Range("F5").Value = 23.92252713947907
$ws.Range("G5").Value = 28.10913579848449
$ws.Range("H5").Value = 14.241512645278
$ws.Range("K5").Value = 7.606865358795231
$ws.Range("L5").Value = 10.65537145869188
$ws.Range("N5").Value = 18.69250058826163
$ws.Range("O5").Value = 21.56975621989253
$ws.Range("B6").Value = 11.48351453837834
$ws.Range("C6").Value = 11.6748702510059
$ws.Range("D6").Value = 4.869612859102154
$ws.Range("F6").Value = 23.92438824710524
$ws.Range("G6").Value = 28.11205836307234
$ws.Range("H6").Value = 14.24345388355409
$ws.Range("K6").Value = 7.598340018300691
$ws.Range("L6").Value = 10.65448242657417
$ws.Range("N6").Value = 18.69499931504703
$ws.Range("O6").Value = 21.5728637552556
$ws.Range("B7").Value = 11.56491490573751
$ws.Range("C7").Value = 11.66233586881889
$ws.Range("D7").Value = 4.891514042947116
$ws.Range("F7").Value = 23.91198955013199
$ws.Range("G7").Value = 28.09250177011763
$ws.Range("H7").Value = 14.23016046303424
$ws.Range("K7").Value = 7.657404684070535
$ws.Range("L7").Value = 10.66084450266132
$ws.Range("N7").Value = 18.677800758899
$ws.Range("O7").Value = 21.55166396063387
$ws.Range("B8").Value = 11.92187006833693
$ws.Range("C8").Value = 11.60979418070695
$ws.Range("D8").Value = 4.985495737175433
$ws.Range("F8").Value = 23.87014158875996
$ws.Range("G8").Value = 28.0242874301389
$ws.Range("H8").Value = 14.1759746251961
$ws.Range("K8").Value = 7.914787853520827
$ws.Range("L8").Value = 10.6936199445697
$ws.Range("N8").Value = 18.60549567187356
$ws.Range("O8").Value = 21.46727656906714
$ws.Range("B9").Value = 12.60711325431105
$ws.Range("C9").Value = 11.51665300415287
$ws.Range("D9").Value = 5.159544950824628
$ws.Range("F9").Value = 23.83438295265411
$ws.Range("G9").Value = 27.95583549672394
$ws.Range("H9").Value = 14.08578988207284
$ws.Range("K9").Value = 8.403171293112964
$ws.Range("L9").Value = 10.77445821681179
$ws.Range("N9").Value = 18.47652476177111
$ws.Range("O9").Value = 21.33481178964564
$ws.Range("B10").Value = 13.09606078108124
$ws.Range("C10").Value = 11.45419368319747
$ws.Range("D10").Value = 5.280606523522687
$ws.Range("F10").Value = 23.83635386311323
$ws.Range("G10").Value = 27.94557233770465
$ws.Range("H10").Value = 14.02931975389757
$ws.Range("K10").Value = 8.793241573035131
$ws.Range("L10").Value = 10.84361287706236
$ws.Range("N10").Value = 18.3895092433621
$ws.Range("O10").Value = 21.25769306734801
$ws.Range("B11").Value = 13.31429175211682
$ws.Range("C11").Value = 11.42706331949991
$ws.Range("D11").Value = 5.334090168071288
$ws.Range("F11").Value = 23.84337225197233
$ws.Range("G11").Value = 27.94961774137384
$ws.Range("H11").Value = 14.00575335630825
$ws.Range("K11").Value = 8.98933286746764
$ws.Range("L11").Value = 10.87711149262621
$ws.Range("N11").Value = 18.3515880973649
$ws.Range("O11").Value = 21.22701076504741
$ws.Range("B12").Value = 13.39624719309979
$ws.Range("C12").Value = 11.41697326207001
$ws.Range("D12").Value = 5.35410591849111
$ws.Range("F12").Value = 23.84690807838353
$ws.Range("G12").Value = 27.95240262728475
$ws.Range("H12").Value = 13.99713439054483
$ws.Range("K12").Value = 9.062264005330206
$ws.Range("L12").Value = 10.89008226017226
$ws.Range("N12").Value = 18.33746625192564
$ws.Range("O12").Value = 21.21602580309803
$ws.Range("B13").Value = 13.37862826117658
$ws.Range("C13").Value = 11.41913818290419
$ws.Range("D13").Value = 5.349805863076343
$ws.Range("F13").Value = 23.84610756320216
$ws.Range("G13").Value = 27.9517471399956
$ws.Range("H13").Value = 13.99897706961159
$ws.Range("K13").Value = 9.046616365047292
$ws.Range("L13").Value = 10.88727620047962
$ws.Range("N13").Value = 18.34049707188583
$ws.Range("O13").Value = 21.21836340773764
$ws.Range("B14").Value = 13.32104842316392
$ws.Range("C14").Value = 11.42622952904576
$ws.Range("D14").Value = 5.33574168263182
$ws.Range("F14").Value = 23.84364556296767
$ws.Range("G14").Value = 27.94982174886411
$ws.Range("H14").Value = 14.00503815406994
$ws.Range("K14").Value = 8.99535967360357
$ws.Range("L14").Value = 10.87817292816082
$ws.Range("N14").Value = 18.35042151916035
$ws.Range("O14").Value = 21.22609431804048
$ws.Range("B15").Value = 13.28568775169797
$ws.Range("C15").Value = 11.43059707196327
$ws.Range("D15").Value = 5.327095797098244
$ws.Range("F15").Value = 23.84225179565627
$ws.Range("G15").Value = 27.94880554271913
$ws.Range("H15").Value = 14.00879047721756
$ws.Range("K15").Value = 8.963790083964696
$ws.Range("L15").Value = 10.87263386318616
$ws.Range("N15").Value = 18.35653150336011
$ws.Range("O15").Value = 21.2309122897582
$ws.Range("B16").Value = 13.08170743609571
$ws.Range("C16").Value = 11.45599245645335
$ws.Range("D16").Value = 5.277078508770029
$ws.Range("F16").Value = 23.83601821045743
$ws.Range("G16").Value = 27.94548333017798
$ws.Range("H16").Value = 14.03090257271809
$ws.Range("K16").Value = 8.780242932072087
$ws.Range("L16").Value = 10.84146405692499
$ws.Range("N16").Value = 18.39202084896212
$ws.Range("O16").Value = 21.25978685440881
$ws.Range("B17").Value = 12.95543855774956
$ws.Range("C17").Value = 11.47189963990363
$ws.Range("D17").Value = 5.245981273681365
$ws.Range("F17").Value = 23.83376072048202
$ws.Range("G17").Value = 27.94567746558342
$ws.Range("H17").Value = 14.04501112955241
$ws.Range("K17").Value = 8.665318099677268
$ws.Range("L17").Value = 10.82285945502147
$ws.Range("N17").Value = 18.41421752990809
$ws.Range("O17").Value = 21.2786280962973
$ws.Range("B18").Value = 12.88242175445412
$ws.Range("C18").Value = 11.48116981335181
$ws.Range("D18").Value = 5.227946018442115
$ws.Range("F18").Value = 23.83303866318253
$ws.Range("G18").Value = 27.94660942654653
$ws.Range("H18").Value = 14.05332574772673
$ws.Range("K18").Value = 8.598373106459555
$ws.Range("L18").Value = 10.81235106896726
$ws.Range("N18").Value = 18.42714103025942
$ws.Range("O18").Value = 21.28987911111779
$ws.Range("B19").Value = 12.85763505081627
$ws.Range("C19").Value = 11.48432930461075
$ws.Range("D19").Value = 5.221814291000153
$ws.Range("F19").Value = 23.8328932376821
$ws.Range("G19").Value = 27.94706584367996
$ws.Range("H19").Value = 14.05617524302536
$ws.Range("K19").Value = 8.580282828825467
$ws.Range("L19").Value = 10.80882639940742
$ws.Range("N19").Value = 18.43154362535751
$ws.Range("O19").Value = 21.29375958523996
$ws.Range("B20").Value = 12.96892109327486
$ws.Range("C20").Value = 11.47019379857709
$ws.Range("D20").Value = 5.249307109735241
$ws.Range("F20").Value = 23.83394139339513
$ws.Range("G20").Value = 27.94557189864513
$ws.Range("H20").Value = 14.04348857657957
$ws.Range("K20").Value = 8.677639577614221
$ws.Range("L20").Value = 10.82482007959903
$ws.Range("N20").Value = 18.41183845971465
$ws.Range("O20").Value = 21.27657955472398
$ws.Range("B21").Value = 13.33798019201367
$ws.Range("C21").Value = 11.42414165188025
$ws.Range("D21").Value = 5.339879184357097
$ws.Range("F21").Value = 23.84434490143032
$ws.Range("G21").Value = 27.9503532841364
$ws.Range("H21").Value = 14.00324958621794
$ws.Range("K21").Value = 9.010451193053203
$ws.Range("L21").Value = 10.88083909352988
$ws.Range("N21").Value = 18.34750001533617
$ws.Range("O21").Value = 21.22380635375952
$ws.Range("B22").Value = 13.57516219715157
$ws.Range("C22").Value = 11.39511382590873
$ws.Range("D22").Value = 5.39768598780613
$ws.Range("F22").Value = 23.85626098770997
$ws.Range("G22").Value = 27.96078062971446
$ws.Range("H22").Value = 13.978729514272
$ws.Range("K22").Value = 9.220231047808467
$ws.Range("L22").Value = 10.91911138389715
$ws.Range("N22").Value = 18.30683831017544
$ws.Range("O22").Value = 21.1930104982215
$ws.Range("B23").Value = 13.44896661196857
$ws.Range("C23").Value = 11.410508897059
$ws.Range("D23").Value = 5.366963203321045
$ws.Range("F23").Value = 23.84943387472318
$ws.Range("G23").Value = 27.9545475280078
$ws.Range("H23").Value = 13.9916536272647
$ws.Range("K23").Value = 9.10898490938345
$ws.Range("L23").Value = 10.89853546366644
$ws.Range("N23").Value = 18.32841363293301
$ws.Range("O23").Value = 21.20910843616787
$ws.Range("B24").Value = 12.96282695074488
$ws.Range("C24").Value = 11.4709646202557
$ws.Range("D24").Value = 5.247803987714128
$ws.Range("F24").Value = 23.83385791743651
$ws.Range("G24").Value = 27.94561707014122
$ws.Range("H24").Value = 14.04417628912995
$ws.Range("K24").Value = 8.672071754620852
$ws.Range("L24").Value = 10.82393309614883
$ws.Range("N24").Value = 18.41291353160943
$ws.Range("O24").Value = 21.27750439550479
$ws.Range("B25").Value = 12.42390557374316
$ws.Range("C25").Value = 11.54079723388359
$ws.Range("D25").Value = 5.11361710918189
$ws.Range("F25").Value = 23.83909520359495
$ws.Range("G25").Value = 27.96733508233049
$ws.Range("H25").Value = 14.1084675138716
$ws.Range("K25").Value = 8.273206052238745
$ws.Range("L25").Value = 10.75084860471718
$ws.Range("N25").Value = 18.51005027751868
$ws.Range("O25").Value = 21.36710412331245
